# P-122 nouvelle structure simple pour l'algo MFCC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New data rows 21 (DCT init) and 22 (DCT 20 -> 20 coeffs).
#    Shared-string insertion order matters (matches the target sharedStrings
#    table ordering), so write C22 before C21.
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = "DCT 20 - > 20 coeffs"
$ws.Range("D22").Value = 20834

$ws.Range("C21").Value = "DCT init"
$ws.Range("D21").Value = 324781

# ---------------------------------------------------------------------------
# 2) New header row 6 ("run" / "init" over F/G) and extra header cell G7.
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = "run"
$ws.Range("G6").Value = "init"
$ws.Range("G7").Value = $ws.Range("F7").Value2

# ---------------------------------------------------------------------------
# 3) Extend the "temps (ms)" column formula (E) as one shared formula across
#    E13:E31 (covers the two new rows + the blank rows down to the new SUM
#    position), then extend the "% cycle" column formula (F) the same way.
# ---------------------------------------------------------------------------
$ws.Range("E13:E31").Formula = "=D13/225000000*1000"
$ws.Range("F13:F31").Formula = "=E13/10"

# The old row-8/9 "% cycle" values are removed (kept blank, same style).
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()

# Rows 12, 16 and 21 move their "% cycle" value from column F to column G
# (independent formulas, not part of the F shared-formula group).
$ws.Range("F12").ClearContents()
$ws.Range("G12").Formula = "=E12/10"

$ws.Range("F16").ClearContents()
$ws.Range("G16").Formula = "=E16/10"

$ws.Range("F21").ClearContents()
$ws.Range("G21").Formula = "=E21/10"

# ---------------------------------------------------------------------------
# 4) Clear the filler/spacer rows 23-31 (formulas were only needed to seed
#    the shared-formula group above) and give them the right number formats.
# ---------------------------------------------------------------------------
$ws.Range("E23:E31").ClearContents()
$ws.Range("F23:F31").ClearContents()

$ws.Range("E23:E31").NumberFormat = "0.00000"
$ws.Range("E23:E31").HorizontalAlignment = -4131
$ws.Range("F23:F31").NumberFormat = "0.00%"
$ws.Range("F23:F31").HorizontalAlignment = -4131

# Row 26 only keeps the E cell - no F cell at all there.
$ws.Range("F26").Clear()

# The running total moves from F32 down to F25.
$ws.Range("F32").Clear()
$ws.Range("F25").Formula = "=SUM(F8:F31)"

# ---------------------------------------------------------------------------
# 5) Column widths: column G gets the same "narrow number" width as column F.
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 11.9

# ---------------------------------------------------------------------------
# 6) View state: active cell / selection.
# ---------------------------------------------------------------------------
$ws.Range("G25").Select()
